$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 327 (pushes existing rows 327..377 down to 328..378,
# inheriting formatting such as the date style on column D from the row above).
$ws.Rows(327).Insert()

# Populate the newly inserted row 327 with the new record.
$ws.Range("A327").Value = 7
$ws.Range("B327").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C327").Value = "Ñuble"
$ws.Range("D327").Value = [DateTime]::ParseExact("2021-10-05", "yyyy-MM-dd", $null)
$ws.Range("E327").Value = 16
$ws.Range("F327").Value = 100112033
$ws.Range("G327").Value = "Lechuga"
$ws.Range("H327").Value = "Conconina(o)"
$ws.Range("I327").Value = "Segunda"
$ws.Range("J327").Value = 120
$ws.Range("K327").Value = 5000
$ws.Range("L327").Value = 5500
$ws.Range("M327").Value = 5250
$ws.Range("N327").Value = "`$/caja 12 unidades"
$ws.Range("O327").Value = "Provincia de Diguillín"
$ws.Range("P327").Value = 438
$ws.Range("Q327").Value = 12
$ws.Range("R327").Value = "Hortaliza"
